$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D:D").Insert()

# Copy number/date formatting from column E (the old column D, now shifted) into the new column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write the updated financial data (new column D + revised figures in the shifted columns)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43100
$ws.Range("F7").Value = 42735
$ws.Range("G7").Value = 42369
$ws.Range("H7").Value = 42004
$ws.Range("I7").Value = 41639
$ws.Range("J7").Value = 41274
$ws.Range("K7").Value = 40908

$ws.Range("D8").Value = 14527000
$ws.Range("E8").Value = 12681000
$ws.Range("F8").Value = 10743000
$ws.Range("G8").Value = 9224000
$ws.Range("H8").Value = 8442000
$ws.Range("I8").Value = 6793300
$ws.Range("J8").Value = 5261000
$ws.Range("K8").Value = 4355600

$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = 242000
$ws.Range("F9").Value = 415000
$ws.Range("G9").Value = 632200
$ws.Range("H9").Value = 857800
$ws.Range("I9").Value = 1077400
$ws.Range("J9").Value = 1177300
$ws.Range("K9").Value = 1275700

$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = 12439000
$ws.Range("F10").Value = 10328000
$ws.Range("G10").Value = 8591800
$ws.Range("H10").Value = 7584100
$ws.Range("I10").Value = 5715900
$ws.Range("J10").Value = 4083700
$ws.Range("K10").Value = 3079900


$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"

$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0

$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 35000
$ws.Range("F14").Value = 1004000
$ws.Range("G14").Value = "NA"
$ws.Range("H14").Value = "NA"
$ws.Range("I14").Value = "NA"
$ws.Range("J14").Value = "NA"
$ws.Range("K14").Value = "NA"

$ws.Range("D15").Value = 426000
$ws.Range("E15").Value = 363000
$ws.Range("F15").Value = 309000
$ws.Range("G15").Value = 272500
$ws.Range("H15").Value = 207800
$ws.Range("I15").Value = 118000
$ws.Range("J15").Value = 65100
$ws.Range("K15").Value = 53800


$ws.Range("D17").Value = 9186000
$ws.Range("E17").Value = 8151000
$ws.Range("F17").Value = 7900000
$ws.Range("G17").Value = 5965100
$ws.Range("H17").Value = 5368700
$ws.Range("I17").Value = 4380900
$ws.Range("J17").Value = 3431200
$ws.Range("K17").Value = 2956700

$ws.Range("D18").Value = 5341000
$ws.Range("E18").Value = 4530000
$ws.Range("F18").Value = 2843000
$ws.Range("G18").Value = 3258900
$ws.Range("H18").Value = 3073300
$ws.Range("I18").Value = 2412400
$ws.Range("J18").Value = 1829800
$ws.Range("K18").Value = 1398900


$ws.Range("D20").Value = -237000
$ws.Range("E20").Value = 123000
$ws.Range("F20").Value = 78000
$ws.Range("G20").Value = 29600
$ws.Range("H20").Value = 4500
$ws.Range("I20").Value = -32600
$ws.Range("J20").Value = -5900
$ws.Range("K20").Value = 600

$ws.Range("D21").Value = 5530000
$ws.Range("E21").Value = 5015800
$ws.Range("F21").Value = 3230100
$ws.Range("G21").Value = 3561000
$ws.Range("H21").Value = 3285600
$ws.Range("I21").Value = 2497800
$ws.Range("J21").Value = 1889100
$ws.Range("K21").Value = 1453300

$ws.Range("D22").Value = 269000
$ws.Range("E22").Value = 254000
$ws.Range("F22").Value = 208000
$ws.Range("G22").Value = 160200
$ws.Range("H22").Value = 88400
$ws.Range("I22").Value = 83300
$ws.Range("J22").Value = 62100
$ws.Range("K22").Value = 31700

$ws.Range("D23").Value = 4835000
$ws.Range("E23").Value = 4399000
$ws.Range("F23").Value = 2713000
$ws.Range("G23").Value = 3128300
$ws.Range("H23").Value = 2989400
$ws.Range("I23").Value = 2296500
$ws.Range("J23").Value = 1761900
$ws.Range("K23").Value = 1367800

$ws.Range("D24").Value = 883000
$ws.Range("E24").Value = 675000
$ws.Range("F24").Value = 578000
$ws.Range("G24").Value = 577000
$ws.Range("H24").Value = 567700
$ws.Range("I24").Value = 403700
$ws.Range("J24").Value = 337800
$ws.Range("K24").Value = 308700

$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0

$ws.Range("D26").Value = 3952000
$ws.Range("E26").Value = 3724000
$ws.Range("F26").Value = 2135000
$ws.Range("G26").Value = 2551400
$ws.Range("H26").Value = 2421800
$ws.Range("I26").Value = 1892800
$ws.Range("J26").Value = 1424000
$ws.Range("K26").Value = 1059100

$ws.Range("D27").Value = 3952000
$ws.Range("E27").Value = 3724000
$ws.Range("F27").Value = 2135000
$ws.Range("G27").Value = 2551400
$ws.Range("H27").Value = 2421800
$ws.Range("I27").Value = 1892700
$ws.Range("J27").Value = 1419600
$ws.Range("K27").Value = 1056400

$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0

$ws.Range("D29").Value = 46000
$ws.Range("E29").Value = -1383000
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "NA"
$ws.Range("H29").Value = "NA"
$ws.Range("I29").Value = "NA"
$ws.Range("J29").Value = "NA"
$ws.Range("K29").Value = "NA"

$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0

$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0

$ws.Range("D32").Value = 237000
$ws.Range("E32").Value = -123000
$ws.Range("F32").Value = -78000
$ws.Range("G32").Value = -29600
$ws.Range("H32").Value = -4500
$ws.Range("I32").Value = 32600
$ws.Range("J32").Value = 5900
$ws.Range("K32").Value = -600

$ws.Range("D33").Value = 3998000
$ws.Range("E33").Value = 2341000
$ws.Range("F33").Value = 2135000
$ws.Range("G33").Value = 2551400
$ws.Range("H33").Value = 2421800
$ws.Range("I33").Value = 1892700
$ws.Range("J33").Value = 1419600
$ws.Range("K33").Value = 1056400

$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0

$ws.Range("D35").Value = 3998000
$ws.Range("E35").Value = 2341000
$ws.Range("F35").Value = 2135000
$ws.Range("G35").Value = 2551400
$ws.Range("H35").Value = 2421800
$ws.Range("I35").Value = 1892700
$ws.Range("J35").Value = 1419600
$ws.Range("K35").Value = 1056400

$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43100
$ws.Range("F38").Value = 42735
$ws.Range("G38").Value = 42369
$ws.Range("H38").Value = 42004
$ws.Range("I38").Value = 41639
$ws.Range("J38").Value = 41274
$ws.Range("K38").Value = 40908



$ws.Range("D41").Value = 2624000
$ws.Range("E41").Value = 2541600
$ws.Range("F41").Value = 2081100
$ws.Range("G41").Value = 1477300
$ws.Range("H41").Value = 3148700
$ws.Range("I41").Value = 1290000
$ws.Range("J41").Value = 1536300
$ws.Range("K41").Value = 632800

$ws.Range("D42").Value = 3660000
$ws.Range("E42").Value = 4859900
$ws.Range("F42").Value = 2218900
$ws.Range("G42").Value = 1171200
$ws.Range("H42").Value = 1142200
$ws.Range("I42").Value = 5462700
$ws.Range("J42").Value = 3646800
$ws.Range("K42").Value = 2024800

$ws.Range("D43").Value = 1523000
$ws.Range("E43").Value = 1217800
$ws.Range("F43").Value = 860100
$ws.Range("G43").Value = 645200
$ws.Range("H43").Value = 643900
$ws.Range("I43").Value = 536000
$ws.Range("J43").Value = 367500
$ws.Range("K43").Value = 264500

$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0

$ws.Range("D45").Value = 600000
$ws.Range("E45").Value = 415500
$ws.Range("F45").Value = 241400
$ws.Range("G45").Value = 259600
$ws.Range("H45").Value = 332600
$ws.Range("I45").Value = 192300
$ws.Range("J45").Value = 131700
$ws.Range("K45").Value = 144700

$ws.Range("D46").Value = 8407000
$ws.Range("E46").Value = 9034800
$ws.Range("F46").Value = 5401500
$ws.Range("G46").Value = 3553200
$ws.Range("H46").Value = 5113600
$ws.Range("I46").Value = 7480900
$ws.Range("J46").Value = 5682400
$ws.Range("K46").Value = 3066800

$ws.Range("D47").Value = 8408000
$ws.Range("E47").Value = 10421600
$ws.Range("F47").Value = 9591100
$ws.Range("G47").Value = 7931400
$ws.Range("H47").Value = 3755700
$ws.Range("I47").Value = "NA"
$ws.Range("J47").Value = "NA"
$ws.Range("K47").Value = "NA"

$ws.Range("D48").Value = 656000
$ws.Range("E48").Value = 480100
$ws.Range("F48").Value = 347000
$ws.Range("G48").Value = 274800
$ws.Range("H48").Value = 199000
$ws.Range("I48").Value = 135100
$ws.Range("J48").Value = 89300
$ws.Range("K48").Value = 64300

$ws.Range("D49").Value = 5035000
$ws.Range("E49").Value = 4914500
$ws.Range("F49").Value = 4390800
$ws.Range("G49").Value = 5542500
$ws.Range("H49").Value = 5661200
$ws.Range("I49").Value = 2787900
$ws.Range("J49").Value = 730800
$ws.Range("K49").Value = 704900

$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0

$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0

$ws.Range("D52").Value = 181000
$ws.Range("E52").Value = 600300
$ws.Range("F52").Value = 108600
$ws.Range("G52").Value = 118700
$ws.Range("H52").Value = 41500
$ws.Range("I52").Value = 40600
$ws.Range("J52").Value = 67300
$ws.Range("K52").Value = 134600

$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0

$ws.Range("D54").Value = 22687000
$ws.Range("E54").Value = 25451300
$ws.Range("F54").Value = 19839000
$ws.Range("G54").Value = 17420600
$ws.Range("H54").Value = 14771000
$ws.Range("I54").Value = 10444500
$ws.Range("J54").Value = 6569700
$ws.Range("K54").Value = 3970700



$ws.Range("D57").Value = 1134000
$ws.Range("E57").Value = 667500
$ws.Range("F57").Value = 419100
$ws.Range("G57").Value = 322800
$ws.Range("H57").Value = 281500
$ws.Range("I57").Value = 247300
$ws.Range("J57").Value = 184600
$ws.Range("K57").Value = 146900

$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 710900
$ws.Range("F58").Value = 967700
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 37100
$ws.Range("I58").Value = 151900
$ws.Range("J58").Value = 520300
$ws.Range("K58").Value = 497600

$ws.Range("D59").Value = 2421000
$ws.Range("E59").Value = 2119400
$ws.Range("F59").Value = 1471800
$ws.Range("G59").Value = 1116500
$ws.Range("H59").Value = 1060100
$ws.Range("I59").Value = 982500
$ws.Range("J59").Value = 756700
$ws.Range("K59").Value = 461300

$ws.Range("D60").Value = 3555000
$ws.Range("E60").Value = 3497900
$ws.Range("F60").Value = 2858700
$ws.Range("G60").Value = 1439300
$ws.Range("H60").Value = 1378700
$ws.Range("I60").Value = 1381700
$ws.Range("J60").Value = 1461700
$ws.Range("K60").Value = 1105800

$ws.Range("D61").Value = 8649000
$ws.Range("E61").Value = 8812800
$ws.Range("F61").Value = 6199100
$ws.Range("G61").Value = 6158400
$ws.Range("H61").Value = 3824200
$ws.Range("I61").Value = 1750600
$ws.Range("J61").Value = 936700
$ws.Range("K61").Value = 77400

$ws.Range("D62").Value = 1698000
$ws.Range("E62").Value = 1880000
$ws.Range("F62").Value = 961100
$ws.Range("G62").Value = 1027400
$ws.Range("H62").Value = 1001400
$ws.Range("I62").Value = 402400
$ws.Range("J62").Value = 114100
$ws.Range("K62").Value = 86200

$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0

$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0

$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0

$ws.Range("D66").Value = 13902000
$ws.Range("E66").Value = 14190700
$ws.Range("F66").Value = 10018800
$ws.Range("G66").Value = 8625100
$ws.Range("H66").Value = 6204300
$ws.Range("I66").Value = 3534700
$ws.Range("J66").Value = 2672800
$ws.Range("K66").Value = 1396400


$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0

$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0

$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0

$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0

$ws.Range("D72").Value = 18367000
$ws.Range("E72").Value = 13938900
$ws.Range("F72").Value = 11326900
$ws.Range("G72").Value = 9191900
$ws.Range("H72").Value = 6640500
$ws.Range("I72").Value = 4218800
$ws.Range("J72").Value = 2368600
$ws.Range("K72").Value = 1033700

$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0

$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0

$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0

$ws.Range("D76").Value = 8785000
$ws.Range("E76").Value = 11260600
$ws.Range("F76").Value = 9820100
$ws.Range("G76").Value = 8795500
$ws.Range("H76").Value = 8566700
$ws.Range("I76").Value = 6909700
$ws.Range("J76").Value = 3897000
$ws.Range("K76").Value = 2574300

$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0

$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43100
$ws.Range("F80").Value = 42735
$ws.Range("G80").Value = 42369
$ws.Range("H80").Value = 42004
$ws.Range("I80").Value = 41639
$ws.Range("J80").Value = 41274
$ws.Range("K80").Value = 40908

$ws.Range("D81").Value = 3998000
$ws.Range("E81").Value = 2341000
$ws.Range("F81").Value = 2135000
$ws.Range("G81").Value = 2551400
$ws.Range("H81").Value = 2421800
$ws.Range("I81").Value = 1892700
$ws.Range("J81").Value = 1419600
$ws.Range("K81").Value = 1056400


$ws.Range("D83").Value = 426000
$ws.Range("E83").Value = 362800
$ws.Range("F83").Value = 309100
$ws.Range("G83").Value = 272500
$ws.Range("H83").Value = 207800
$ws.Range("I83").Value = 118000
$ws.Range("J83").Value = 65100
$ws.Range("K83").Value = 53800

$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0

$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0

$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0

$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0

$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0

$ws.Range("D89").Value = 5338000
$ws.Range("E89").Value = 4662000
$ws.Range("F89").Value = 3924700
$ws.Range("G89").Value = 3102200
$ws.Range("H89").Value = 2914400
$ws.Range("I89").Value = 2301400
$ws.Range("J89").Value = 1785800
$ws.Range("K89").Value = 1341800


$ws.Range("D91").Value = -442000
$ws.Range("E91").Value = -287800
$ws.Range("F91").Value = -219900
$ws.Range("G91").Value = -173900
$ws.Range("H91").Value = -131500
$ws.Range("I91").Value = -84400
$ws.Range("J91").Value = -55200
$ws.Range("K91").Value = -46800

$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0

$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0

$ws.Range("D94").Value = 2215000
$ws.Range("E94").Value = -4202000
$ws.Range("F94").Value = -3333300
$ws.Range("G94").Value = -3894500
$ws.Range("H94").Value = -2357900
$ws.Range("I94").Value = -2162300
$ws.Range("J94").Value = -1562700
$ws.Range("K94").Value = -904800


$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0

$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0

$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0

$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0

$ws.Range("D100").Value = -7431000
$ws.Range("E100").Value = -78700
$ws.Range("F100").Value = 57700
$ws.Range("G100").Value = -730000
$ws.Range("H100").Value = 1429000
$ws.Range("I100").Value = -403500
$ws.Range("J100").Value = 668900
$ws.Range("K100").Value = -151000

$ws.Range("D101").Value = -40000
$ws.Range("E101").Value = 100000
$ws.Range("F101").Value = -45200
$ws.Range("G101").Value = -149100
$ws.Range("H101").Value = -136500
$ws.Range("I101").Value = 18000
$ws.Range("J101").Value = 11600
$ws.Range("K101").Value = -12100

$ws.Range("D102").Value = 82000
$ws.Range("E102").Value = 481300
$ws.Range("F102").Value = 603900
$ws.Range("G102").Value = -1671400
$ws.Range("H102").Value = 1849000
$ws.Range("I102").Value = -246400
$ws.Range("J102").Value = 903500
$ws.Range("K102").Value = 273900
